$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.248.27"
$ws.Range("E2").Value = "  +2.44%  "

$ws.Range("D3").Value = "1.906.09"
$ws.Range("E3").Value = "  +2.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.62"
$ws.Range("E5").Value = "  +1.20%  "

$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4656"
$ws.Range("E7").Value = "  +1.61%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3965"
$ws.Range("E8").Value = "  +2.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.71"
$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07968"
$ws.Range("E10").Value = "  +1.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +3.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.35"
$ws.Range("E12").Value = "  +2.11%  "

$ws.Range("D13").Value = "1.914.86"
$ws.Range("E13").Value = "  +3.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.145"
$ws.Range("E14").Value = "  +2.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.784"
$ws.Range("E15").Value = "  +1.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06951"
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.92"
$ws.Range("E17").Value = "  +0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001012"
$ws.Range("E19").Value = "  +1.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.19"
$ws.Range("E20").Value = "  +2.51%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").Value = "29.263.73"
$ws.Range("E22").Value = "  +2.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.360"
$ws.Range("E23").Value = "  +2.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").Value = "  +0.81%  "

$ws.Range("D25").Value = "2.138.45"
$ws.Range("E25").Value = "  +3.09%  "

$ws.Range("E26").Value = "  -2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.67"
$ws.Range("E27").Value = "  +2.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.52"
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.913"
$ws.Range("E29").Value = "  +2.49%  "

$ws.Range("E30").Value = "  +1.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.57"
$ws.Range("E31").Value = "  +0.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09443"
$ws.Range("E32").Value = "  +1.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9279"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("E34").Value = "  +1.80%  "

$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("E36").Value = "  -1.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05864"
$ws.Range("E37").Value = "  +1.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.172"
$ws.Range("E38").Value = "  +2.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02108"
$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.997"
$ws.Range("E40").Value = "  +4.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5762"
$ws.Range("E41").Value = "  +2.81%  "

$ws.Range("E42").Value = "  +2.02%  "

$ws.Range("E43").Value = "  +3.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.08"
$ws.Range("E44").Value = "  +3.00%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.238"
$ws.Range("E45").Value = "  +4.39%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5435"
$ws.Range("E46").Value = "  +3.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07094"
$ws.Range("E47").Value = "  -1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.887"
$ws.Range("E48").Value = "  +3.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.595"
$ws.Range("E49").Value = "  +7.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.97"
$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.070"
$ws.Range("E51").Value = "  -5.72%  "
